# Generate Report for Handback
# Updates the generated-report timestamps on the Overview, zh-cn and de-de
# sheets of the handback-status workbook.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 87f810d5-e8df-4903-b300-3720d7e05e96.md row.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-30 18:56:50"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 87f810d5-e8df-4903-b300-3720d7e05e96 handback file row.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-30 18:56:45"
$wsZhCn.Range("K4").Value = "2016-08-30 18:57:14"

# de-de sheet: same "Latest HO Xliff Generate Date" value as the Overview
# sheet, plus its own "Correspond Handback DateTime".
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-30 18:56:50"
$wsDeDe.Range("K4").Value = "2016-08-30 18:57:22"
